# Apply updated auto-eval results to top_3_configurations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: meta_llama3_instruct_70B, top_p 0.5, new bge-large embedding ---
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = "('bge-large', 'bge-large-en-v1.5', 1024)"
$ws.Range("F2").Value = 3.8
$ws.Range("G2").Value = 80
$ws.Range("J2").Value = 3.8
$ws.Range("K2").Value = 80
$ws.Range("N2").Value = "/Users/yonghuizhu/imperial/contoso-chat-backend/eval/auto_eval/meta_llama3_instruct_70B_top0.5_emb('bge-large', 'bge-large-en-v1.5', 1024)_originaltemplate.ipynb"
$ws.Range("O2").Value = "2024_06_19_155426chat_eval_run"

# --- Row 3: Phi_3_mini_4k_instruct, top_p 0.5, new bge-large embedding ---
$ws.Range("A3").Value = "Phi_3_mini_4k_instruct"
$ws.Range("B3").Value = 4000
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = "['bge-large', 'bge-large-en-v1.5', 1024]"
$ws.Range("F3").Value = 4.2
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3.4
$ws.Range("K3").Value = 60
$ws.Range("L3").Value = 1.8
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "/Users/yonghuizhu/imperial/contoso-chat-backend/eval/auto_eval/Phi_3_mini_4k_instruct_top0.5_emb['bge-large', 'bge-large-en-v1.5', 1024]_originaltemplate.ipynb"
$ws.Range("O3").Value = "2024_06_19_152518chat_eval_run"
$ws.Range("P3").Value = 18.08

# --- Row 4: meta_llama3_instruct_70B, top_p 0.9, text-embedding-ada-002 (previous row2 config moved down) ---
$ws.Range("A4").Value = "meta_llama3_instruct_70B"
$ws.Range("B4").Value = 8000
$ws.Range("C4").Value = 0.9
$ws.Range("F4").Value = 4.2
$ws.Range("G4").Value = 100
$ws.Range("J4").Value = 3.4
$ws.Range("N4").Value = "/Users/yonghuizhu/imperial/contoso-chat-backend/eval/auto_eval/meta_llama3_instruct_70B_top0.9_embtext-embedding-ada-002_originaltemplate.ipynb"
$ws.Range("O4").Value = "2024_06_19_131434chat_eval_run"
$ws.Range("P4").Value = 18.08

$wb.Save()
